$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: PHP, updated value
$ws.Range("B2").Value = 2065

# Row 3: now Python (was node.js), updated value
$ws.Range("A3").Value = "Python"
$ws.Range("B3").Value = 2112

# Row 4: now node.js (was Python), updated value
$ws.Range("A4").Value = "node.js"
$ws.Range("B4").Value = 2267

# Update selection
$ws.Range("A11").Select()
